$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.040.91"
$ws.Range("D3").Value = "1.648.84"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.14"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5193"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2628"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06312"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.36"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07657"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.578"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.57%  "
$ws.Range("D13").Value = "1.644.13"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").Value = "1.876.58"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5579"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "0.0₅8121"
$ws.Range("E16").Value = "  +1.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.15"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").Value = "26.028.48"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.605"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.50"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "191.00"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.899"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.75"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1181"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.183"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.83"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  +1.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05344"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.269"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.451"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.341"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.551"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.94%  "
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.780"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9425"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5620"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01573"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.872"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").Value = "1.029.93"
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8218"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.54"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.84%  "
$ws.Range("D45").Value = "1.787.08"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("E46").Value = "  +7.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.17"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4314"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.922"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E51").Value = "  -3.61%  "
